# Add "Exp Constant" and "Exp Constant [dB]" columns (AA, AB) to the table,
# mirroring the formatting of the existing header row and filling the
# computed sphere position / experimental-constant values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Range("AA1").Value = "Exp Constant"
$ws.Range("AB1").Value = "Exp Constant [dB]"

# Copy the formatting (bold font, border, centered/top alignment) from the
# existing header cell Z1 onto the two new header cells so they match the
# rest of the header row.
$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)

# --- Data cells (rows 2-32) ---
# Same constant values for every row.
$ws.Range("AA2:AA32").Value = 385250961.9682089
$ws.Range("AB2:AB32").Value = 85.85743731821252
